# Migración a la nueva interfaz:
#  - Renombra el encabezado J1 ("Vencimiento" -> "Vencimiento_Formateada")
#  - Agrega una nueva columna K ("Estado") marcando cada fila como "Vencido"
#  - Resalta la columna "Haber" (H): rojo si tiene saldo (<> 0), verde si esta en 0
#  - (replica el comportamiento observado en el commit original) C7 tambien
#    queda resaltada en rojo junto con el resto de las celdas "Vencido"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$RED = 255        # BGR-packed 0x0000FF -> fgColor FF0000
$GREEN = 582476    # BGR-packed 0x08E34C -> fgColor 4CE308

# 1) Encabezados: renombrar J1 y agregar K1 copiando el formato del resto de encabezados
$ws.Range("J1").Value = "Vencimiento_Formateada"

$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Estado"

# 2) Filas de datos: nueva columna "Estado" + resaltado de la columna "Haber"
$lastRow = 8
for ($r = 2; $r -le $lastRow; $r++) {
    $haber = $ws.Cells.Item($r, 8).Value2

    if ($haber -ne 0) {
        $ws.Cells.Item($r, 8).Interior.Color = $RED
    } else {
        $ws.Cells.Item($r, 8).Interior.Color = $GREEN
    }

    $estadoCell = $ws.Cells.Item($r, 11)
    $estadoCell.Value = "Vencido"
    $estadoCell.Interior.Color = $RED
}

# La celda C7 tambien quedo resaltada en rojo al migrar la planilla
$ws.Cells.Item(7, 3).Interior.Color = $RED

Write-Output "migracion completa"
